{"js": "// Load every paragraph in the document body so we can locate the two\n// paragraphs referenced by the edit: the one to highlight, and the one\n// after which a new bullet should be inserted.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\n// 1) Highlight the \"Visualize tower recharges\" paragraph in yellow.\n//    Setting highlightColor on the paragraph (rather than on a text range)\n//    applies it both to the run text and to the paragraph mark, matching\n//    how the sibling \"Tower recharges\" bullet above it is already marked.\nconst rechargeParagraph = paragraphs.items.find(\n  (p) => p.text.trim() === \"Visualize tower recharges\"\n);\nif (!rechargeParagraph) {\n  throw new Error('Could not find paragraph \"Visualize tower recharges\".');\n}\nrechargeParagraph.font.highlightColor = \"yellow\";\n\n// 2) Add a new bullet (\"Make proper icons for towers with resource cost\")\n//    right after \"Different enemies (different amounts of health for now)\".\n//    insertParagraph inherits the anchor paragraph's list formatting\n//    (ListParagraph style, ilvl 0, numId 2), so no extra formatting calls\n//    are needed.\nconst enemiesParagraph = paragraphs.items.find(\n  (p) => p.text.trim() === \"Different enemies (different amounts of health for now)\"\n);\nif (!enemiesParagraph) {\n  throw new Error(\n    'Could not find paragraph \"Different enemies (different amounts of health for now)\".'\n  );\n}\nenemiesParagraph.insertParagraph(\n  \"Make proper icons for towers with resource cost\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Highlight the \"Visualize tower recharges\" paragraph in yellow.\n#    Applying the highlight through Range.Font (rather than Range directly)\n#    marks both the run text and the paragraph mark, matching how the\n#    sibling \"Tower recharges\" bullet above it is already highlighted.\n$rechargeParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r\", \"`a\") -eq \"Visualize tower recharges\") {\n        $rechargeParagraph = $p\n        break\n    }\n}\nif ($null -eq $rechargeParagraph) {\n    throw 'Could not find paragraph \"Visualize tower recharges\".'\n}\n$rechargeParagraph.Range.Font.HighlightColorIndex = 7\n\n# 2) Add a new bullet (\"Make proper icons for towers with resource cost\")\n#    right after \"Different enemies (different amounts of health for now)\".\n#    InsertParagraphAfter inherits the source paragraph's list formatting\n#    (ListParagraph style, ilvl 0, numId 2), so no extra formatting calls\n#    are needed.\n$enemiesIndex = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    if ($p.Range.Text.TrimEnd(\"`r\", \"`a\") -eq \"Different enemies (different amounts of health for now)\") {\n        $enemiesIndex = $i\n        break\n    }\n}\nif ($enemiesIndex -eq -1) {\n    throw 'Could not find paragraph \"Different enemies (different amounts of health for now)\".'\n}\n$enemiesParagraph = $d.Paragraphs($enemiesIndex)\n$enemiesParagraph.Range.InsertParagraphAfter()\n# Re-fetch by index: the newly inserted paragraph lives right after it, and\n# objects fetched before the insertion (e.g. via .Next()) do not reliably\n# reflect subsequent writes.\n$newParagraph = $d.Paragraphs($enemiesIndex + 1)\n$newParagraph.Range.Text = \"Make proper icons for towers with resource cost\"\n"}
